# Apply updated dSF (column F) values as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -9
    14 = -3
    17 = 0
    20 = -7
    22 = -3
    23 = -5
    25 = -1
    31 = 5
    33 = -6
    37 = -5
    46 = 0
    49 = -2
    51 = 4
    57 = -2
    58 = 6
    59 = 2
    60 = -4
    61 = 0
    62 = -2
    64 = -4
    65 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
